$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / inline-string cells ---
# J2: DATE_TYPE_CODE "002" -> "001"
# Force a Text number format before assigning so the leading zeros are not
# stripped by automatic General-number coercion, then clear the format again
# so the cell keeps its original (unstyled) appearance.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()

# N2: REPORT_DATE "2020-06-30 00:00:00" -> "2017-12-31 00:00:00"
$ws.Range("N2").Value = "2017-12-31 00:00:00"

# --- Numeric cells ---
$ws.Range("O2").Value = -44125029.51    # PARENT_NETPROFIT
$ws.Range("P2").Value = 143337667.45    # TOTAL_OPERATE_INCOME
$ws.Range("Q2").Value = 180305197.21    # TOTAL_OPERATE_COST
$ws.Range("S2").Value = 114323781.57    # OPERATE_COST
$ws.Range("T2").Value = 114323781.57    # OPERATE_EXPENSE
$ws.Range("V2").Value = 13852909.83     # SALE_EXPENSE
$ws.Range("W2").Value = 18249371.06     # MANAGE_EXPENSE
$ws.Range("X2").Value = 3515562.64      # FINANCE_EXPENSE
$ws.Range("Y2").Value = -53722201.8     # OPERATE_PROFIT
$ws.Range("Z2").Value = -54081440.17    # TOTAL_PROFIT
$ws.Range("AA2").Value = -9956410.66    # INCOME_TAX
$ws.Range("AG2").Value = 1084943.82     # OPERATE_TAX_ADD
$ws.Range("AS2").Value = -46503129.51   # DEDUCT_PARENT_NETPROFIT
